$d = $word.ActiveDocument

# Title date update
$d.Content.Find.Execute("2025-06-26 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-27 Friday", 2) | Out-Null

# Simple one-to-one cell text replacements (unambiguous across the document)
$d.Content.Find.Execute("911÷5=182, 1", $true, $false, $false, $false, $false, $true, 1, $false, "540÷6=90, 0", 2) | Out-Null
$d.Content.Find.Execute("513÷3=171, 0", $true, $false, $false, $false, $false, $true, 1, $false, "711÷3=237, 0", 2) | Out-Null
$d.Content.Find.Execute("225÷6=37, 3", $true, $false, $false, $false, $false, $true, 1, $false, "918÷9=102, 0", 2) | Out-Null
$d.Content.Find.Execute("356÷7=50, 6", $true, $false, $false, $false, $false, $true, 1, $false, "549÷9=61, 0", 2) | Out-Null
$d.Content.Find.Execute("297÷6=49, 3", $true, $false, $false, $false, $false, $true, 1, $false, "407÷3=135, 2", 2) | Out-Null
$d.Content.Find.Execute("751÷5=150, 1", $true, $false, $false, $false, $false, $true, 1, $false, "180÷7=25, 5", 2) | Out-Null
$d.Content.Find.Execute("146÷4=36, 2", $true, $false, $false, $false, $false, $true, 1, $false, "891÷3=297, 0", 2) | Out-Null
$d.Content.Find.Execute("178÷6=29, 4", $true, $false, $false, $false, $false, $true, 1, $false, "498÷6=83, 0", 2) | Out-Null
$d.Content.Find.Execute("764÷4=191, 0", $true, $false, $false, $false, $false, $true, 1, $false, "687÷2=343, 1", 2) | Out-Null
$d.Content.Find.Execute("427÷3=142, 1", $true, $false, $false, $false, $false, $true, 1, $false, "845÷7=120, 5", 2) | Out-Null
$d.Content.Find.Execute("725÷9=80, 5", $true, $false, $false, $false, $false, $true, 1, $false, "591÷6=98, 3", 2) | Out-Null
$d.Content.Find.Execute("974÷2=487, 0", $true, $false, $false, $false, $false, $true, 1, $false, "603÷4=150, 3", 2) | Out-Null
$d.Content.Find.Execute("947÷6=157, 5", $true, $false, $false, $false, $false, $true, 1, $false, "688÷9=76, 4", 2) | Out-Null
$d.Content.Find.Execute("974÷6=162, 2", $true, $false, $false, $false, $false, $true, 1, $false, "899÷8=112, 3", 2) | Out-Null
$d.Content.Find.Execute("898÷2=449, 0", $true, $false, $false, $false, $false, $true, 1, $false, "969÷8=121, 1", 2) | Out-Null
$d.Content.Find.Execute("645÷5=129, 0", $true, $false, $false, $false, $false, $true, 1, $false, "424÷5=84, 4", 2) | Out-Null
$d.Content.Find.Execute("563÷3=187, 2", $true, $false, $false, $false, $false, $true, 1, $false, "467÷5=93, 2", 2) | Out-Null
$d.Content.Find.Execute("615÷7=87, 6", $true, $false, $false, $false, $false, $true, 1, $false, "381÷6=63, 3", 2) | Out-Null
$d.Content.Find.Execute("797÷5=159, 2", $true, $false, $false, $false, $false, $true, 1, $false, "837÷9=93, 0", 2) | Out-Null
$d.Content.Find.Execute("131÷2=65, 1", $true, $false, $false, $false, $false, $true, 1, $false, "343÷4=85, 3", 2) | Out-Null
$d.Content.Find.Execute("759÷2=379, 1", $true, $false, $false, $false, $false, $true, 1, $false, "559÷2=279, 1", 2) | Out-Null
$d.Content.Find.Execute("280÷3=93, 1", $true, $false, $false, $false, $false, $true, 1, $false, "759÷9=84, 3", 2) | Out-Null

# Special row: cell text values shift (925/4->222/5, 913/8 cell now holds old 263/4 text, 263/4 cell gets new text)
# Handled via direct table-cell addressing to avoid ambiguity with duplicate "263÷4=65, 3" text
$t = $d.Tables.Item(1)
$t.Cell(9, 3).Range.Text = "222÷5=44, 2"
$t.Cell(9, 4).Range.Text = "263÷4=65, 3"
$t.Cell(9, 5).Range.Text = "571÷9=63, 4"

